$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Update B6 value (Mass nacelle) from 800000 to 40996
$ws.Range("B6").Value = 40996

# Move the active selection from B5 to B6
$ws.Range("B6").Select()

$wb.Application.Calculate()
